$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows: Instance(col B as text), Clientes(C), Depositos(D), Vehiculos(E), Objective(F), Gap(G), First Sol(H), Time(I)
# A column = a repeating run index (1..4 or 1..2) per instance group

$rows = @(
    @{A=1; B="10x2_1"; C=10; D=2; E=2; F=6311.973347009111;  G=0; H=0.02194023132324219; I=0.04687309265136719},
    @{A=2; B="10x2_1"; C=10; D=2; E=2; F=117.9733470091115;  G=0; H=0.02094268798828125; I=0.04687309265136719},
    @{A=3; B="10x2_1"; C=10; D=2; E=2; F=2522.026652990888;  G=0; H=0.02090644836425781; I=0.04883193969726562},
    @{A=4; B="10x2_1"; C=10; D=2; E=2; F=0.1014301971078715; G=0; H=0.03789901733398438; I=0.09075546264648438},
    @{A=1; B="10x2_2"; C=10; D=2; E=2; F=5860.527052656848;  G=0; H=0.02718353271484375; I=0.0640411376953125},
    @{A=2; B="10x2_2"; C=10; D=2; E=2; F=198.5270526568471;  G=0; H=0.02792167663574219; I=0.06083488464355469},
    @{A=3; B="10x2_2"; C=10; D=2; E=2; F=2319.472947343153;  G=0; H=0.02393341064453125; I=0.05684661865234375},
    @{A=4; B="10x2_2"; C=10; D=2; E=2; F=0.1689638276296211; G=0; H=0.02493476867675781; I=0.06682205200195312},
    @{A=1; B="10x2_3"; C=10; D=2; E=2; F=5925.849472106609;  G=0; H=0.01296234130859375; I=0.04986381530761719},
    @{A=2; B="10x2_3"; C=10; D=2; E=2; F=87.84947210660982;  G=0; H=0.01396369934082031; I=0.06754684448242188}
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r,1).Value = $row.A
    $ws.Cells.Item($r,2).Value = $row.B
    $ws.Cells.Item($r,3).Value = $row.C
    $ws.Cells.Item($r,4).Value = $row.D
    $ws.Cells.Item($r,5).Value = $row.E
    $ws.Cells.Item($r,6).Value = $row.F
    $ws.Cells.Item($r,7).Value = $row.G
    $ws.Cells.Item($r,8).Value = $row.H
    $ws.Cells.Item($r,9).Value = $row.I
    $r++
}
